# Unit Test Design update.
#
# The "Clase" column (2nd column) of the first table's 4 data rows all
# contain the run text "Restaurant". The edit renames these to
# "RestaurantTest" (the actual class name under test), and since Word's
# proofing flags the resulting camel-case word as a possible spelling
# issue, it is wrapped in <w:proofErr w:type="spellStart"/> ... <w:proofErr
# w:type="spellEnd"/>, matching the convention already used elsewhere in
# this document for other camel-case identifiers.
#
# Row 2 keeps "Restaurant" as its own run and simply appends a new run
# "Test" (so the cell ends up with two runs: "Restaurant" + "Test").
# Rows 3-5 collapse straight to a single run "RestaurantTest".

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 2: split into two runs, "Restaurant" + "Test".
$cell = $t.Cell(2, 2)
$null = $cell.Range.InsertXML("<w:p $wNs><w:proofErr w:type='spellStart'/><w:r><w:t>Restaurant</w:t></w:r><w:r><w:t>Test</w:t></w:r><w:proofErr w:type='spellEnd'/></w:p>")

# Rows 3-5: single run "RestaurantTest".
for ($r = 3; $r -le 5; $r++) {
    $cell = $t.Cell($r, 2)
    $null = $cell.Range.InsertXML("<w:p $wNs><w:proofErr w:type='spellStart'/><w:r><w:t>RestaurantTest</w:t></w:r><w:proofErr w:type='spellEnd'/></w:p>")
}
